$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.102636098861694
$ws.Range("B1").Value = 1.960917592048645
$ws.Range("C1").Value = 4.397706508636475
$ws.Range("D1").Value = 0.2356551587581635
$ws.Range("E1").Value = 0.2710447609424591
